# Commit: "Type3 ready? fixing details"
# - Fix the typo/accent in the shared string used by B1 ("Acuífero" -> "Acuifero")
# - Update the sheet's saved selection/active cell (was G12) to F6

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Correct the header text in B1 (removes the accent on the "i")
$ws.Range("B1").Value = "Acuifero"

# Move/record the active selection on the sheet to F6
$ws.Range("F6").Select()
